# Generate Report for Handoff
# Updates the localization-status report: a new source file UUID/hash
# replaces the old one, handoff/handback timestamps are refreshed, and
# the now-unset "Latest Target File" / "Latest Handback File" columns
# are cleared back out for the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "df4e5c99-e057-4b1b-bb78-837172eff53c"
$newGuid = "0753162b-28a5-4a19-8fac-4a6900ecbfff"
$oldHash = "b6451b5773a334ee600a3d527af2cc6c67f8512f"
$newHash = "2d169d0d2d5b7acd8617e5b2c3e413bc85f30213"

$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37395a8d858c0ddaf12e43f43f4e14dede6a6c0f/e2e/$oldGuid.md"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "e2e\$newGuid.md"
$ws1.Range("G2").Value = "2016-09-03 01:04:21"

# Re-point the B2 hyperlink at the new display text while keeping the
# same underlying target / relationship.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $sourceUrl, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newGuid.md")

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-03 01:04:16"
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

# I2 no longer carries a hyperlink/target file, drop its hyperlink and
# restore plain "Normal" formatting; keep the A2 hyperlink alive.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $sourceUrl, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md")
$ws2.Range("I2").Style = "Normal"

$ws2.Columns.Item(9).ColumnWidth = 17.8333333333333
$ws2.Columns.Item(10).ColumnWidth = 20.8333333333333

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-03 01:04:21"
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $sourceUrl, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md")
$ws3.Range("I2").Style = "Normal"

$ws3.Columns.Item(9).ColumnWidth = 17.8333333333333
$ws3.Columns.Item(10).ColumnWidth = 20.8333333333333
